$d = $word.ActiveDocument

# The resume's opening line used to be a Heading1 "Matthew Mikolajczyk" paragraph
# wrapping a "matthew-mikolajczyk" bookmark. Thea's demo replaces it with a plain
# FirstParagraph-styled banner line and drops the bookmark entirely (the remaining
# section bookmarks simply shift down to fill the gap, which Word does on its own).

# 1. Drop the now-unused "matthew-mikolajczyk" bookmark.
$d.Bookmarks("matthew-mikolajczyk").Delete()

# 2. Swap the heading text for the demo banner text (Range.Text keeps the run's
#    xml:space="preserve" intact, unlike Find/Replace).
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Text = "HELLO THEA Matthew Mikolajczyk ==================="

# 3. Re-style the paragraph from Heading1 down to FirstParagraph.
$titlePara.Style = "FirstParagraph"

Write-Output $d.Paragraphs(1).Range.Text
